# feat: add 2022-Q4 data
#
# Before: 3 sheets -> 总计, 2022-Q3, 2022-Q2
# After:  4 sheets -> 总计, 2022-Q4, 2022-Q3, 2022-Q2
#
# The new "2022-Q4" sheet carries fresh fund numbers; "2022-Q3" keeps
# its original data (it's really a duplicate of the old "2022-Q3"
# sheet, since that slot gets renamed/overwritten into "2022-Q4");
# "2022-Q2" is untouched. The "总计" (totals) sheet gains a new
# summary row for 2022-Q2 and its existing rows shift down one slot.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

# Duplicate "2022-Q3" and place the copy right after it -- the copy
# keeps the old data and becomes the new "2022-Q3" sheet, freeing the
# original to be turned into "2022-Q4".
$q3.Copy($null, $q3)
$q3copy = $wb.Worksheets.Item("2022-Q3 (2)")

# Rename the original sheet to "2022-Q4" first so the duplicate can
# take over the "2022-Q3" name without a collision.
$q3.Name = "2022-Q4"
$q3copy.Name = "2022-Q3"

# Update the fund metrics on the new "2022-Q4" sheet. These columns
# are stored as text in the workbook, so force a text number format
# before writing the numeric-looking strings (otherwise Excel would
# silently convert them to numbers).
$q4Metrics = $q3.Range("D2:G3")
$q4Metrics.NumberFormat = "@"

$q3.Range("D2").Value = "1.79"
$q3.Range("E2").Value = "88.58"
$q3.Range("F2").Value = "3.93"
$q3.Range("G2").Value = "0.0703"
$q3.Range("H2").Value = 8

$q3.Range("D3").Value = "1.79"
$q3.Range("E3").Value = "88.58"
$q3.Range("F3").Value = "3.93"
$q3.Range("G3").Value = "0.0703"
$q3.Range("H3").Value = 8

# Update the "总计" summary sheet: row 2 is now Q4, row 3 is now Q3
# (same 0.14 value as before), and a brand-new row 4 carries the Q2
# totals that used to live in row 3. Copy row 3's formatting down to
# the new row 4 first so the new cells match the existing style.
$total.Range("A3").Copy($total.Range("A4"))

$total.Range("B2").Value = "2022-Q4"
$total.Range("B3").Value = "2022-Q3"
$total.Range("D3").Value = 0.14

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.16
